# Update forecast data: insert a new earliest row (2007/2008 forecast) and
# shift the remaining rows by one, refreshing all the y_0_forecast /
# y_1_forecast values across the whole series (bugfix to evaluation and
# simulated rt_data for components).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(39400,2007,7.226520411029047,2008,13.12477391005418),
    @(39583,2008,7.424215059809214,2009,14.69926045795804),
    @(39765,2008,4.268860212333636,2009,13.08276537368067),
    @(39948,2009,-8.992252553594248,2010,-19.76480035196674),
    @(40130,2009,-7.266312015249776,2010,12.31225042954258),
    @(40310,2010,7.007132997505194,2011,26.6762552377682),
    @(40494,2010,6.958243460951929,2011,12.04357532583245),
    @(40676,2011,8.86644397614711,2012,5.982869268853186),
    @(40862,2011,9.469137444079934,2012,10.6992064911972),
    @(41044,2012,3.001306214623578,2013,0.08911682035466217),
    @(41228,2012,3.358206407534947,2013,4.390489499870109),
    @(41409,2013,-2.904769335987201,2014,-8.13908495190001),
    @(41592,2013,0.3081076735359067,2014,3.50243535103556),
    @(41774,2014,6.148460028297587,2015,9.131012060398724),
    @(41957,2014,3.901355411819707,2015,6.143002545701282),
    @(42137,2015,5.171596082708629,2016,5.986046065677453),
    @(42321,2015,5.331683351557981,2016,4.555278923792572),
    @(42503,2016,4.337699953939178,2017,5.703761500281579),
    @(42689,2016,3.254758369308375,2017,0.9515943257393689),
    @(42867,2017,4.073887526082065,2018,1.602827009554897),
    @(43053,2017,5.246209615995667,2018,4.251116704684899),
    @(43145,2018,7.340964210079881,2019,8.197760099691198),
    @(43235,2018,2.281540236993274,2019,-4.308894244053674),
    @(43326,2018,5.15263050486201,2019,6.518301903862911),
    @(43418,2018,4.86255966374296,2019,4.992093705734701),
    @(43510,2019,3.497157880977619,2020,2.970919592955878),
    @(43600,2019,3.436682959168125,2020,2.82953744009995),
    @(43691,2019,2.727571990295941,2020,-0.7915059299106297),
    @(43783,2019,2.764740011159428,2020,0.3611963426346065),
    @(43875,2020,2.93113923573054,2021,5.185727774795068),
    @(43966,2020,-2.305533699949835,2021,-6.248031846400004),
    @(44068,2020,-7.578477024949737,2021,17.08516853885251),
    @(44159,2020,-7.260793671746435,2021,2.387971016884638),
    @(44251,2021,-0.8212867190516282,2022,-3.88653222228782),
    @(44341,2021,5.121832664816339,2022,12.25634856540583),
    @(44432,2021,4.379227219808146,2022,3.504204983279191),
    @(44525,2021,4.097586525396268,2022,3.347989317130651),
    @(44617,2022,5.102173676573241,2023,-5.324806173859886),
    @(44706,2022,8.333485306093348,2023,7.042341419899389),
    @(44798,2022,7.397318165265498,2023,1.541970075766486),
    @(44890,2022,7.824284864703746,2023,-2.122471977790918),
    @(44981,2023,-0.9123477982208139,2024,-5.598178600215808),
    @(45071,2023,0.08405665459807476,2024,0.4484505192704713),
    @(45163,2023,-0.3046246622258053,2024,1.814369137627603),
    @(45254,2023,-1.24502235313334,2024,-4.678511595261359),
    @(45345,2024,-3.378144228902036,2025,-0.366247000203368),
    @(45436,2024,-2.055826728150212,2025,3.376972582720295),
    @(45534,2024,-2.567041707495976,2025,-0.6241557099098238),
    @(45618,2024,-1.735114423676209,2025,1.758477003221981),
    @(45713,2025,1.830872485486124,2026,-0.1271413384857256),
    @(45800,2025,1.378024997308636,2026,0.4446401485209472),
    @(45891,2025,2.64031107104763,2026,2.021190955273178),
)

$startRow = 2
$endRow = $startRow + $data.Count - 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# Make sure the date-formatted style used by column A is applied to the
# newly added row as well (copy style from the row above it, values were
# already written above so only formats get overwritten here).
$ws.Range("A" + ($endRow - 1)).Copy()
$ws.Range("A$endRow").PasteSpecial(-4122)
$excel.CutCopyMode = 0
